$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('!!Compartment')
$ws.Range('A1').Value = '!!!ObjTables schema=''SBtab'' objTablesVersion=''0.0.8'' date=''2020-03-09 23:58:44'''

$ws = $wb.Worksheets.Item('!!Compartment')
$ws.Range('A2').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Compartment'' name=''Compartment'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Compound')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Compound'' name=''Compound'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Definition')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Definition'' name=''Definition'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Enzyme')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Enzyme'' name=''Enzyme'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!FbcObjective')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''FbcObjective'' name=''FbcObjective'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Gene')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Gene'' name=''Gene'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Layout')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Layout'' name=''Layout'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Measurement')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Measurement'' name=''Measurement'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!PbConfig')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''PbConfig'' name=''PbConfig'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Position')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Position'' name=''Position'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Protein')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Protein'' name=''Protein'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Quantity')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Quantity'' name=''Quantity'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!QuantityInfo')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''QuantityInfo'' name=''QuantityInfo'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!QuantityMatrix')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''QuantityMatrix'' name=''QuantityMatrix'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Reaction')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Reaction'' name=''Reaction'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!ReactionStoichiometry')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''ReactionStoichiometry'' name=''ReactionStoichiometry'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Regulator')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Regulator'' name=''Regulator'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Relation')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Relation'' name=''Relation'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!Relationship')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''Relationship'' name=''Relationship'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!SparseMatrix')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''SparseMatrix'' name=''SparseMatrix'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!SparseMatrixColumn')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''SparseMatrixColumn'' name=''SparseMatrixColumn'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!SparseMatrixOrdered')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''SparseMatrixOrdered'' name=''SparseMatrixOrdered'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!SparseMatrixRow')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''SparseMatrixRow'' name=''SparseMatrixRow'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!StoichiometricMatrix')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''StoichiometricMatrix'' name=''StoichiometricMatrix'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!rxnconContingencyList')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''rxnconContingencyList'' name=''rxnconContingencyList'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''

$ws = $wb.Worksheets.Item('!!rxnconReactionList')
$ws.Range('A1').Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' id=''rxnconReactionList'' name=''rxnconReactionList'' date=''2020-03-09 23:58:44'' objTablesVersion=''0.0.8'''
